$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column D (Template_name) values to be prefixed with the Study_Type (column C)
$ws.Range("D2").Value = "Clinical_search-strategy-template.xlsx"
$ws.Range("D3").Value = "Economic_search-strategy-template.xlsx"
$ws.Range("D4").Value = "Quality of life_search-strategy-template.xlsx"
$ws.Range("D5").Value = "Real-world Evidence_search-strategy-template.xlsx"

# Update the selected range on the sheet to D2:D5
$ws.Range("D2:D5").Select()
